# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Column G ("K") is recomputed; write the new values for rows 2..69.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(
    0,3,0,1,2,0,2,0,4,0,
    2,3,2,1,2,1,1,3,2,1,
    2,0,0,1,1,2,0,0,1,1,
    3,1,1,3,4,0,1,0,1,1,
    1,1,1,1,1,0,0,1,1,2,
    0,0,1,0,1,2,4,1,1,0,
    3,2,2,0,1,2,2,1
)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
